$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "scenario2" marker in column A for rows 8-10 (matching the
# existing pattern used for rows 2-5 / row 7).
$ws.Range("A8").Value = "scenario2"
$ws.Range("A9").Value = "scenario2"
$ws.Range("A10").Value = "scenario2"

# Update the active selection to match the new state (B12).
$ws.Range("B12").Select()
